# Insert a new weekly price record as row 221 in the daily-logic subset
# sheet, pushing the existing rows 221-274 down to 222-275 (dimension
# grows from R274 to R275).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 221:274 down by one, creating a blank row 221.
$ws.Rows("221:221").Insert()

# Populate the newly inserted row with the new market record.
$ws.Range("A221").Value = 4
$ws.Range("B221").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C221").Value = "Los Lagos"
$ws.Range("D221").Value = 44722
$ws.Range("E221").Value = 10
$ws.Range("F221").Value = 100112043
$ws.Range("G221").Value = "Pepino ensalada"
$ws.Range("H221").Value = "Sin especificar"
$ws.Range("I221").Value = "Primera"
$ws.Range("J221").Value = 350
$ws.Range("K221").Value = 23000
$ws.Range("L221").Value = 23000
$ws.Range("M221").Value = 23000
$ws.Range("N221").Value = "$/caja 60 unidades"
$ws.Range("O221").Value = "Región de Arica y Parinacota"
$ws.Range("P221").Value = 383
$ws.Range("Q221").Value = 60
$ws.Range("R221").Value = "Hortaliza"
